$wb = $excel.ActiveWorkbook

# "Input-files" sheet: document the new punch-motion input file.
$ws = $wb.Worksheets.Item("Input-files")

# Row 7 already has "pellet_v7_20180706.i" in column A; add its description.
$ws.Range("B7").Value = "penalty based boundary condition implemented. It allows to turn on & off temp bc based on spatial location."

# New row 8: new input file name + description.
$ws.Range("A8").Value = "pellet_v8_20180612.i"
$ws.Range("B8").Value = "punch motion implemented"

# Make "Input-files" the active sheet/tab and leave the selection on the
# first empty row below the data that was just entered.
$ws.Activate()
$ws.Range("A9").Select()
